$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Testing the Nesto App login page"
$ws.Range("C6").Value = '5.Verify that the URL is "dashboard"'

$ws.Range("A7").Value = "TC_02"
$ws.Range("B7").Value = "Invalid Login Test"
$ws.Range("C7").Value = '1. Open URL "http://localhost:8080/login"'
$ws.Range("C8").Value = '2.Type "shana@gmail.com" into "//input[@placeholder=''Enter your email'']"'
$ws.Range("C9").Value = '3.Type "wrongpass123" into "//input[@placeholder=''Enter your password'']"'
$ws.Range("C10").Value = '4.Click on the "Sign In" button "//button[@type=''submit'']"'
$ws.Range("C11").Value = '5.Verify URL is "error"'
$ws.Range("C12").Value = '6.Verify text "Invalid Email or Password" at "//div[contains(@class, ''error-msg'')]"'

$ws.Range("A13").Value = "TC_03"
$ws.Range("B13").Value = "Check Signup Link"
$ws.Range("C13").Value = '1.Open URL "http://localhost:8080/login"'
$ws.Range("C14").Value = '2.Click on "Register Link" at "//a[@href=''/signup'']"'
$ws.Range("C15").Value = '3.Verify URL is "signup"'

$ws.Range("C23").Select()

$ws.PageSetup.Orientation = 1

